$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.277433666666667
$ws.Cells.Item(2, 8).Value = 3.832301
$ws.Cells.Item(2, 9).Value = 0.01786062203930835
$ws.Cells.Item(2, 10).Value = 0.01786062203930835
$ws.Cells.Item(2, 13).Value = 33.51218733333334
$ws.Cells.Item(2, 14).Value = 100.536562
$ws.Cells.Item(2, 15).Value = 0.5521050876757374
$ws.Cells.Item(2, 16).Value = 0.5521050876757374
$ws.Cells.Item(2, 17).Value = 42.80959634324022
$ws.Cells.Item(2, 18).Value = 385.286367089162
$ws.Cells.Item(2, 19).Value = 0.009860940296955547
$ws.Cells.Item(2, 20).Value = 0.009860940296955546
$ws.Cells.Item(3, 7).Value = 1.277433666666667
$ws.Cells.Item(3, 8).Value = 3.832301
$ws.Cells.Item(3, 9).Value = 0.01786062203930835
$ws.Cells.Item(3, 10).Value = 0.01786062203930835
$ws.Cells.Item(3, 14).Value = 54.5272
$ws.Cells.Item(3, 15).Value = 0.2994407600362589
$ws.Cells.Item(3, 16).Value = 0.299440760036259
$ws.Cells.Item(3, 17).Value = 23.21829367635555
$ws.Cells.Item(3, 18).Value = 208.9646430872
$ws.Cells.Item(3, 19).Value = 0.00534819823817085
$ws.Cells.Item(3, 20).Value = 0.00534819823817085
$ws.Cells.Item(4, 7).Value = 1.277433666666667
$ws.Cells.Item(4, 8).Value = 3.832301
$ws.Cells.Item(4, 9).Value = 0.01786062203930835
$ws.Cells.Item(4, 10).Value = 0.01786062203930835
$ws.Cells.Item(4, 13).Value = 2.975281333333333
$ws.Cells.Item(4, 14).Value = 8.925844
$ws.Cells.Item(4, 15).Value = 0.04901703207436071
$ws.Cells.Item(4, 16).Value = 0.04901703207436071
$ws.Cells.Item(4, 17).Value = 3.800724543004888
$ws.Cells.Item(4, 18).Value = 34.20652088704399
$ws.Cells.Item(4, 19).Value = 0.0008754746833688115
$ws.Cells.Item(4, 20).Value = 0.0008754746833688113
$ws.Cells.Item(5, 7).Value = 1.277433666666667
$ws.Cells.Item(5, 8).Value = 3.832301
$ws.Cells.Item(5, 9).Value = 0.01786062203930835
$ws.Cells.Item(5, 10).Value = 0.01786062203930835
$ws.Cells.Item(5, 13).Value = 6.035726666666666
$ws.Cells.Item(5, 14).Value = 18.10718
$ws.Cells.Item(5, 15).Value = 0.09943712021364286
$ws.Cells.Item(5, 16).Value = 0.09943712021364286
$ws.Cells.Item(5, 17).Value = 7.710240446797776
$ws.Cells.Item(5, 18).Value = 69.39216402117999
$ws.Cells.Item(5, 19).Value = 0.001776008820813144
$ws.Cells.Item(5, 20).Value = 0.001776008820813144
$ws.Cells.Item(6, 7).Value = 45.44725166666667
$ws.Cells.Item(6, 9).Value = 0.6354272679079697
$ws.Cells.Item(6, 10).Value = 0.6354272679079697
$ws.Cells.Item(6, 13).Value = 33.51218733333334
$ws.Cells.Item(6, 14).Value = 100.536562
$ws.Cells.Item(6, 15).Value = 0.5521050876757374
$ws.Cells.Item(6, 16).Value = 0.5521050876757374
$ws.Cells.Item(6, 17).Value = 1523.036811638479
$ws.Cells.Item(6, 18).Value = 13707.33130474631
$ws.Cells.Item(6, 19).Value = 0.3508226274598839
$ws.Cells.Item(6, 20).Value = 0.3508226274598839
$ws.Cells.Item(7, 7).Value = 45.44725166666667
$ws.Cells.Item(7, 9).Value = 0.6354272679079697
$ws.Cells.Item(7, 10).Value = 0.6354272679079697
$ws.Cells.Item(7, 14).Value = 54.5272
$ws.Cells.Item(7, 15).Value = 0.2994407600362589
$ws.Cells.Item(7, 16).Value = 0.299440760036259
$ws.Cells.Item(7, 17).Value = 826.0371270262223
$ws.Cells.Item(7, 18).Value = 7434.334143236
$ws.Cells.Item(7, 19).Value = 0.1902728240501259
$ws.Cells.Item(7, 20).Value = 0.190272824050126
$ws.Cells.Item(8, 7).Value = 45.44725166666667
$ws.Cells.Item(8, 9).Value = 0.6354272679079697
$ws.Cells.Item(8, 10).Value = 0.6354272679079697
$ws.Cells.Item(8, 13).Value = 2.975281333333333
$ws.Cells.Item(8, 14).Value = 8.925844
$ws.Cells.Item(8, 15).Value = 0.04901703207436071
$ws.Cells.Item(8, 16).Value = 0.04901703207436071
$ws.Cells.Item(8, 17).Value = 135.2183595351356
$ws.Cells.Item(8, 18).Value = 1216.96523581622
$ws.Cells.Item(8, 19).Value = 0.03114675877196835
$ws.Cells.Item(8, 20).Value = 0.03114675877196835
$ws.Cells.Item(9, 7).Value = 45.44725166666667
$ws.Cells.Item(9, 9).Value = 0.6354272679079697
$ws.Cells.Item(9, 10).Value = 0.6354272679079697
$ws.Cells.Item(9, 13).Value = 6.035726666666666
$ws.Cells.Item(9, 14).Value = 18.10718
$ws.Cells.Item(9, 15).Value = 0.09943712021364286
$ws.Cells.Item(9, 16).Value = 0.09943712021364286
$ws.Cells.Item(9, 17).Value = 274.3071888112111
$ws.Cells.Item(9, 18).Value = 2468.7646993009
$ws.Cells.Item(9, 19).Value = 0.06318505762599143
$ws.Cells.Item(9, 20).Value = 0.06318505762599143
$ws.Cells.Item(10, 7).Value = 23.96074166666667
$ws.Cells.Item(10, 8).Value = 71.88222500000001
$ws.Cells.Item(10, 9).Value = 0.3350105464235513
$ws.Cells.Item(10, 10).Value = 0.3350105464235513
$ws.Cells.Item(10, 13).Value = 33.51218733333334
$ws.Cells.Item(10, 14).Value = 100.536562
$ws.Cells.Item(10, 15).Value = 0.5521050876757374
$ws.Cells.Item(10, 16).Value = 0.5521050876757374
$ws.Cells.Item(10, 17).Value = 802.976863378939
$ws.Cells.Item(10, 18).Value = 7226.791770410451
$ws.Cells.Item(10, 19).Value = 0.1849610271054715
$ws.Cells.Item(10, 20).Value = 0.1849610271054715
$ws.Cells.Item(11, 7).Value = 23.96074166666667
$ws.Cells.Item(11, 8).Value = 71.88222500000001
$ws.Cells.Item(11, 9).Value = 0.3350105464235513
$ws.Cells.Item(11, 10).Value = 0.3350105464235513
$ws.Cells.Item(11, 14).Value = 54.5272
$ws.Cells.Item(11, 15).Value = 0.2994407600362589
$ws.Cells.Item(11, 16).Value = 0.299440760036259
$ws.Cells.Item(11, 17).Value = 435.5040510022222
$ws.Cells.Item(11, 18).Value = 3919.53645902
$ws.Cells.Item(11, 19).Value = 0.1003158126412306
$ws.Cells.Item(11, 20).Value = 0.1003158126412306
$ws.Cells.Item(12, 7).Value = 23.96074166666667
$ws.Cells.Item(12, 8).Value = 71.88222500000001
$ws.Cells.Item(12, 9).Value = 0.3350105464235513
$ws.Cells.Item(12, 10).Value = 0.3350105464235513
$ws.Cells.Item(12, 13).Value = 2.975281333333333
$ws.Cells.Item(12, 14).Value = 8.925844
$ws.Cells.Item(12, 15).Value = 0.04901703207436071
$ws.Cells.Item(12, 16).Value = 0.04901703207436071
$ws.Cells.Item(12, 17).Value = 71.28994741365555
$ws.Cells.Item(12, 18).Value = 641.6095267229
$ws.Cells.Item(12, 19).Value = 0.01642122269929232
$ws.Cells.Item(12, 20).Value = 0.01642122269929232
$ws.Cells.Item(13, 7).Value = 23.96074166666667
$ws.Cells.Item(13, 8).Value = 71.88222500000001
$ws.Cells.Item(13, 9).Value = 0.3350105464235513
$ws.Cells.Item(13, 10).Value = 0.3350105464235513
$ws.Cells.Item(13, 13).Value = 6.035726666666666
$ws.Cells.Item(13, 14).Value = 18.10718
$ws.Cells.Item(13, 15).Value = 0.09943712021364286
$ws.Cells.Item(13, 16).Value = 0.09943712021364286
$ws.Cells.Item(13, 17).Value = 144.6204874306111
$ws.Cells.Item(13, 18).Value = 1301.5843868755
$ws.Cells.Item(13, 19).Value = 0.03331248397755685
$ws.Cells.Item(13, 20).Value = 0.03331248397755685
$ws.Cells.Item(14, 7).Value = 0.8369233333333334
$ws.Cells.Item(14, 8).Value = 2.51077
$ws.Cells.Item(14, 9).Value = 0.01170156362917063
$ws.Cells.Item(14, 10).Value = 0.01170156362917063
$ws.Cells.Item(14, 13).Value = 33.51218733333334
$ws.Cells.Item(14, 14).Value = 100.536562
$ws.Cells.Item(14, 15).Value = 0.5521050876757374
$ws.Cells.Item(14, 16).Value = 0.5521050876757374
$ws.Cells.Item(14, 17).Value = 28.04713153030445
$ws.Cells.Item(14, 18).Value = 252.42418377274
$ws.Cells.Item(14, 19).Value = 0.006460492813426472
$ws.Cells.Item(14, 20).Value = 0.006460492813426471
$ws.Cells.Item(15, 7).Value = 0.8369233333333334
$ws.Cells.Item(15, 8).Value = 2.51077
$ws.Cells.Item(15, 9).Value = 0.01170156362917063
$ws.Cells.Item(15, 10).Value = 0.01170156362917063
$ws.Cells.Item(15, 14).Value = 54.5272
$ws.Cells.Item(15, 15).Value = 0.2994407600362589
$ws.Cells.Item(15, 16).Value = 0.299440760036259
$ws.Cells.Item(15, 17).Value = 15.21169532711111
$ws.Cells.Item(15, 18).Value = 136.905257944
$ws.Cells.Item(15, 19).Value = 0.003503925106731498
$ws.Cells.Item(15, 20).Value = 0.003503925106731498
$ws.Cells.Item(16, 7).Value = 0.8369233333333334
$ws.Cells.Item(16, 8).Value = 2.51077
$ws.Cells.Item(16, 9).Value = 0.01170156362917063
$ws.Cells.Item(16, 10).Value = 0.01170156362917063
$ws.Cells.Item(16, 13).Value = 2.975281333333333
$ws.Cells.Item(16, 14).Value = 8.925844
$ws.Cells.Item(16, 15).Value = 0.04901703207436071
$ws.Cells.Item(16, 16).Value = 0.04901703207436071
$ws.Cells.Item(16, 17).Value = 2.490082371097778
$ws.Cells.Item(16, 18).Value = 22.41074133988
$ws.Cells.Item(16, 19).Value = 0.0005735759197312296
$ws.Cells.Item(16, 20).Value = 0.0005735759197312295
$ws.Cells.Item(17, 7).Value = 0.8369233333333334
$ws.Cells.Item(17, 8).Value = 2.51077
$ws.Cells.Item(17, 9).Value = 0.01170156362917063
$ws.Cells.Item(17, 10).Value = 0.01170156362917063
$ws.Cells.Item(17, 13).Value = 6.035726666666666
$ws.Cells.Item(17, 14).Value = 18.10718
$ws.Cells.Item(17, 15).Value = 0.09943712021364286
$ws.Cells.Item(17, 16).Value = 0.09943712021364286
$ws.Cells.Item(17, 17).Value = 5.051440480955556
$ws.Cells.Item(17, 18).Value = 45.46296432859999
$ws.Cells.Item(17, 19).Value = 0.001163569789281431
$ws.Cells.Item(17, 20).Value = 0.001163569789281431
